{"js": "// Replace the 100 arithmetic answers in the single table, cell-by-cell,\n// preserving row/column position (row-major order) so that duplicate\n// values before/after the edit don't get cross-clobbered by a naive\n// global find/replace.\nconst oldValues = [\n  \"94-62=32\", \"53+16=69\", \"89-63=26\", \"98-81=17\", \"38+1=39\",\n  \"13+9=22\", \"79-69=10\", \"34+24=58\", \"4+19=23\", \"54+19=73\",\n  \"2+71=73\", \"83-12=71\", \"75-2=73\", \"18+71=89\", \"71+19=90\",\n  \"21+17=38\", \"34+39=73\", \"52-49=3\", \"57+8=65\", \"9+22=31\",\n  \"35+6=41\", \"51+26=77\", \"41-24=17\", \"94-71=23\", \"55-33=22\",\n  \"29+26=55\", \"13+43=56\", \"88-73=15\", \"31+58=89\", \"45+51=96\",\n  \"79-24=55\", \"99-32=67\", \"57-54=3\", \"53+14=67\", \"12+75=87\",\n  \"66+19=85\", \"16+25=41\", \"43+45=88\", \"63-56=7\", \"29+20=49\",\n  \"81-69=12\", \"42+17=59\", \"7+76=83\", \"42-17=25\", \"85-56=29\",\n  \"95-21=74\", \"37-18=19\", \"36-3=33\", \"84-73=11\", \"21+11=32\",\n  \"11-7=4\", \"50+2=52\", \"99-27=72\", \"80-62=18\", \"23-10=13\",\n  \"10+59=69\", \"39+3=42\", \"63+27=90\", \"74-24=50\", \"40+6=46\",\n  \"36+42=78\", \"38+36=74\", \"59-21=38\", \"32+49=81\", \"19-13=6\",\n  \"88-55=33\", \"27+62=89\", \"35-21=14\", \"27+71=98\", \"23+17=40\",\n  \"59-9=50\", \"11-2=9\", \"59-16=43\", \"60+14=74\", \"83-13=70\",\n  \"5+62=67\", \"16-1=15\", \"85-72=13\", \"49+41=90\", \"10+76=86\",\n  \"3+68=71\", \"48+4=52\", \"61-59=2\", \"93-55=38\", \"60-53=7\",\n  \"50+0=50\", \"31+22=53\", \"18-9=9\", \"47+16=63\", \"34+3=37\",\n  \"99-24=75\", \"81-10=71\", \"37+22=59\", \"47+41=88\", \"87-60=27\",\n  \"40+47=87\", \"24+11=35\", \"78+16=94\", \"94-26=68\", \"73-36=37\"\n];\n\nconst newValues = [\n  \"27+71=98\", \"59-30=29\", \"32+4=36\", \"4+46=50\", \"74-73=1\",\n  \"77+9=86\", \"80-18=62\", \"3+90=93\", \"36+11=47\", \"69-46=23\",\n  \"47-15=32\", \"36+56=92\", \"11+59=70\", \"86-13=73\", \"72-36=36\",\n  \"43+5=48\", \"25+68=93\", \"95-24=71\", \"49+43=92\", \"46+34=80\",\n  \"83-67=16\", \"46-18=28\", \"14+60=74\", \"94-46=48\", \"99-81=18\",\n  \"40+26=66\", \"41+42=83\", \"29+42=71\", \"3+48=51\", \"21+75=96\",\n  \"24+8=32\", \"13+33=46\", \"36-24=12\", \"70-61=9\", \"21+52=73\",\n  \"59+37=96\", \"68-40=28\", \"82-55=27\", \"61-37=24\", \"0+57=57\",\n  \"0+89=89\", \"91-73=18\", \"15+80=95\", \"63-32=31\", \"32-0=32\",\n  \"28+4=32\", \"79-3=76\", \"71-62=9\", \"87-69=18\", \"87-79=8\",\n  \"47+42=89\", \"29+61=90\", \"79-6=73\", \"76-73=3\", \"87-33=54\",\n  \"44+53=97\", \"89+3=92\", \"34-25=9\", \"84-64=20\", \"58+25=83\",\n  \"93-6=87\", \"40-13=27\", \"90-73=17\", \"49+1=50\", \"40+7=47\",\n  \"1+52=53\", \"94-32=62\", \"11+83=94\", \"0+35=35\", \"94-11=83\",\n  \"4+47=51\", \"73-42=31\", \"37+23=60\", \"7+9=16\", \"17+64=81\",\n  \"14-13=1\", \"64+21=85\", \"37+23=60\", \"54-16=38\", \"67-41=26\",\n  \"18-11=7\", \"60+36=96\", \"24-14=10\", \"88-45=43\", \"85-21=64\",\n  \"97-25=72\", \"78-65=13\", \"97+1=98\", \"22+45=67\", \"94-62=32\",\n  \"68+23=91\", \"61+31=92\", \"29-13=16\", \"79-66=13\", \"81+6=87\",\n  \"65+5=70\", \"56-43=13\", \"2+5=7\", \"25+31=56\", \"55+12=67\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst rows = table.values.length;\nconst cols = table.values[0].length;\nconst updated = [];\nlet k = 0;\nfor (let r = 0; r < rows; r++) {\n  const row = [];\n  for (let c = 0; c < cols; c++) {\n    const current = table.values[r][c];\n    const expectedOld = oldValues[k];\n    if (current === expectedOld) {\n      row.push(newValues[k]);\n    } else {\n      // Fall back: keep whatever is already there if it doesn't match the\n      // expected \"before\" snapshot (keeps the script idempotent/safe).\n      row.push(current);\n    }\n    k++;\n  }\n  updated.push(row);\n}\n\ntable.values = updated;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic answers in the single table, cell-by-cell,\n# preserving row/column position (row-major order, matching Word's\n# Cell(row, col) 1-based indexing) so duplicate values before/after the\n# edit don't get cross-clobbered by a naive global Find/Replace.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$oldValues = @(\n    \"94-62=32\",\n    \"53+16=69\",\n    \"89-63=26\",\n    \"98-81=17\",\n    \"38+1=39\",\n    \"13+9=22\",\n    \"79-69=10\",\n    \"34+24=58\",\n    \"4+19=23\",\n    \"54+19=73\",\n    \"2+71=73\",\n    \"83-12=71\",\n    \"75-2=73\",\n    \"18+71=89\",\n    \"71+19=90\",\n    \"21+17=38\",\n    \"34+39=73\",\n    \"52-49=3\",\n    \"57+8=65\",\n    \"9+22=31\",\n    \"35+6=41\",\n    \"51+26=77\",\n    \"41-24=17\",\n    \"94-71=23\",\n    \"55-33=22\",\n    \"29+26=55\",\n    \"13+43=56\",\n    \"88-73=15\",\n    \"31+58=89\",\n    \"45+51=96\",\n    \"79-24=55\",\n    \"99-32=67\",\n    \"57-54=3\",\n    \"53+14=67\",\n    \"12+75=87\",\n    \"66+19=85\",\n    \"16+25=41\",\n    \"43+45=88\",\n    \"63-56=7\",\n    \"29+20=49\",\n    \"81-69=12\",\n    \"42+17=59\",\n    \"7+76=83\",\n    \"42-17=25\",\n    \"85-56=29\",\n    \"95-21=74\",\n    \"37-18=19\",\n    \"36-3=33\",\n    \"84-73=11\",\n    \"21+11=32\",\n    \"11-7=4\",\n    \"50+2=52\",\n    \"99-27=72\",\n    \"80-62=18\",\n    \"23-10=13\",\n    \"10+59=69\",\n    \"39+3=42\",\n    \"63+27=90\",\n    \"74-24=50\",\n    \"40+6=46\",\n    \"36+42=78\",\n    \"38+36=74\",\n    \"59-21=38\",\n    \"32+49=81\",\n    \"19-13=6\",\n    \"88-55=33\",\n    \"27+62=89\",\n    \"35-21=14\",\n    \"27+71=98\",\n    \"23+17=40\",\n    \"59-9=50\",\n    \"11-2=9\",\n    \"59-16=43\",\n    \"60+14=74\",\n    \"83-13=70\",\n    \"5+62=67\",\n    \"16-1=15\",\n    \"85-72=13\",\n    \"49+41=90\",\n    \"10+76=86\",\n    \"3+68=71\",\n    \"48+4=52\",\n    \"61-59=2\",\n    \"93-55=38\",\n    \"60-53=7\",\n    \"50+0=50\",\n    \"31+22=53\",\n    \"18-9=9\",\n    \"47+16=63\",\n    \"34+3=37\",\n    \"99-24=75\",\n    \"81-10=71\",\n    \"37+22=59\",\n    \"47+41=88\",\n    \"87-60=27\",\n    \"40+47=87\",\n    \"24+11=35\",\n    \"78+16=94\",\n    \"94-26=68\",\n    \"73-36=37\"\n)\n\n$newValues = @(\n    \"27+71=98\",\n    \"59-30=29\",\n    \"32+4=36\",\n    \"4+46=50\",\n    \"74-73=1\",\n    \"77+9=86\",\n    \"80-18=62\",\n    \"3+90=93\",\n    \"36+11=47\",\n    \"69-46=23\",\n    \"47-15=32\",\n    \"36+56=92\",\n    \"11+59=70\",\n    \"86-13=73\",\n    \"72-36=36\",\n    \"43+5=48\",\n    \"25+68=93\",\n    \"95-24=71\",\n    \"49+43=92\",\n    \"46+34=80\",\n    \"83-67=16\",\n    \"46-18=28\",\n    \"14+60=74\",\n    \"94-46=48\",\n    \"99-81=18\",\n    \"40+26=66\",\n    \"41+42=83\",\n    \"29+42=71\",\n    \"3+48=51\",\n    \"21+75=96\",\n    \"24+8=32\",\n    \"13+33=46\",\n    \"36-24=12\",\n    \"70-61=9\",\n    \"21+52=73\",\n    \"59+37=96\",\n    \"68-40=28\",\n    \"82-55=27\",\n    \"61-37=24\",\n    \"0+57=57\",\n    \"0+89=89\",\n    \"91-73=18\",\n    \"15+80=95\",\n    \"63-32=31\",\n    \"32-0=32\",\n    \"28+4=32\",\n    \"79-3=76\",\n    \"71-62=9\",\n    \"87-69=18\",\n    \"87-79=8\",\n    \"47+42=89\",\n    \"29+61=90\",\n    \"79-6=73\",\n    \"76-73=3\",\n    \"87-33=54\",\n    \"44+53=97\",\n    \"89+3=92\",\n    \"34-25=9\",\n    \"84-64=20\",\n    \"58+25=83\",\n    \"93-6=87\",\n    \"40-13=27\",\n    \"90-73=17\",\n    \"49+1=50\",\n    \"40+7=47\",\n    \"1+52=53\",\n    \"94-32=62\",\n    \"11+83=94\",\n    \"0+35=35\",\n    \"94-11=83\",\n    \"4+47=51\",\n    \"73-42=31\",\n    \"37+23=60\",\n    \"7+9=16\",\n    \"17+64=81\",\n    \"14-13=1\",\n    \"64+21=85\",\n    \"37+23=60\",\n    \"54-16=38\",\n    \"67-41=26\",\n    \"18-11=7\",\n    \"60+36=96\",\n    \"24-14=10\",\n    \"88-45=43\",\n    \"85-21=64\",\n    \"97-25=72\",\n    \"78-65=13\",\n    \"97+1=98\",\n    \"22+45=67\",\n    \"94-62=32\",\n    \"68+23=91\",\n    \"61+31=92\",\n    \"29-13=16\",\n    \"79-66=13\",\n    \"81+6=87\",\n    \"65+5=70\",\n    \"56-43=13\",\n    \"2+5=7\",\n    \"25+31=56\",\n    \"55+12=67\"\n)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n$k = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $expectedOld = $oldValues[$k]\n        $newVal = $newValues[$k]\n\n        # Cell().Range.Text includes a trailing cell-mark (CR + BEL); strip it\n        # before comparing so the lookup is robust.\n        $current = $cell.Range.Text\n        $current = $current.TrimEnd([char]7)\n        $current = $current.TrimEnd([char]13)\n\n        if ($current -eq $expectedOld) {\n            $cell.Range.Text = $newVal\n        }\n        $k++\n    }\n}\n"}
